$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.87751133333333
$ws.Range("H2").Value = 41.632534
$ws.Range("I2").Value = 0.03879702495420557
$ws.Range("J2").Value = 0.03879702495420557
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.275702333333334
$ws.Range("N2").Value = 18.827107
$ws.Range("O2").Value = 0.2560867246427467
$ws.Range("P2").Value = 0.2560867246427466
$ws.Range("Q2").Value = 87.09113025545979
$ws.Range("R2").Value = 783.8201722991381
$ws.Range("S2").Value = 0.009935403046405413
$ws.Range("T2").Value = 0.009935403046405411
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.87751133333333
$ws.Range("H3").Value = 41.632534
$ws.Range("I3").Value = 0.03879702495420557
$ws.Range("J3").Value = 0.03879702495420557
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.25829966666667
$ws.Range("N3").Value = 30.774899
$ws.Range("O3").Value = 0.418600854933333
$ws.Range("P3").Value = 0.418600854933333
$ws.Range("Q3").Value = 142.3596698848962
$ws.Range("R3").Value = 1281.237028964066
$ws.Range("S3").Value = 0.01624046781470031
$ws.Range("T3").Value = 0.01624046781470031
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.87751133333333
$ws.Range("H4").Value = 41.632534
$ws.Range("I4").Value = 0.03879702495420557
$ws.Range("J4").Value = 0.03879702495420557
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.972158333333333
$ws.Range("N4").Value = 23.916475
$ws.Range("O4").Value = 0.3253124204239203
$ws.Range("P4").Value = 0.3253124204239203
$ws.Range("Q4").Value = 110.6337176219611
$ws.Range("R4").Value = 995.70345859765
$ws.Range("S4").Value = 0.01262115409309985
$ws.Range("T4").Value = 0.01262115409309985
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 311.415324
$ws.Range("H5").Value = 934.2459719999999
$ws.Range("I5").Value = 0.8706163379113565
$ws.Range("J5").Value = 0.8706163379113564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.275702333333334
$ws.Range("N5").Value = 18.827107
$ws.Range("O5").Value = 0.2560867246427467
$ws.Range("P5").Value = 0.2560867246427466
$ws.Range("Q5").Value = 1954.349875462556
$ws.Range("R5").Value = 17589.148879163
$ws.Range("S5").Value = 0.2229532863961821
$ws.Range("T5").Value = 0.222953286396182
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 311.415324
$ws.Range("H6").Value = 934.2459719999999
$ws.Range("I6").Value = 0.8706163379113565
$ws.Range("J6").Value = 0.8706163379113564
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.25829966666667
$ws.Range("N6").Value = 30.774899
$ws.Range("O6").Value = 0.418600854933333
$ws.Range("P6").Value = 0.418600854933333
$ws.Range("Q6").Value = 3194.591714384092
$ws.Range("R6").Value = 28751.32542945683
$ws.Range("S6").Value = 0.3644407433686214
$ws.Range("T6").Value = 0.3644407433686213
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 311.415324
$ws.Range("H7").Value = 934.2459719999999
$ws.Range("I7").Value = 0.8706163379113565
$ws.Range("J7").Value = 0.8706163379113564
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.972158333333333
$ws.Range("N7").Value = 23.916475
$ws.Range("O7").Value = 0.3253124204239203
$ws.Range("P7").Value = 0.3253124204239203
$ws.Range("Q7").Value = 2482.6522703543
$ws.Range("R7").Value = 22343.8704331887
$ws.Range("S7").Value = 0.2832223081465531
$ws.Range("T7").Value = 0.283222308146553
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 32.40240933333333
$ws.Range("H8").Value = 97.207228
$ws.Range("I8").Value = 0.09058663713443794
$ws.Range("J8").Value = 0.09058663713443794
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.275702333333334
$ws.Range("N8").Value = 18.827107
$ws.Range("O8").Value = 0.2560867246427467
$ws.Range("P8").Value = 0.2560867246427466
$ws.Range("Q8").Value = 203.3478758588218
$ws.Range("R8").Value = 1830.130882729396
$ws.Range("S8").Value = 0.02319803520015922
$ws.Range("T8").Value = 0.02319803520015921
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 32.40240933333333
$ws.Range("H9").Value = 97.207228
$ws.Range("I9").Value = 0.09058663713443794
$ws.Range("J9").Value = 0.09058663713443794
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.25829966666667
$ws.Range("N9").Value = 30.774899
$ws.Range("O9").Value = 0.418600854933333
$ws.Range("P9").Value = 0.418600854933333
$ws.Range("Q9").Value = 332.3936248633302
$ws.Range("R9").Value = 2991.542623769972
$ws.Range("S9").Value = 0.03791964375001134
$ws.Range("T9").Value = 0.03791964375001133
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.40240933333333
$ws.Range("H10").Value = 97.207228
$ws.Range("I10").Value = 0.09058663713443794
$ws.Range("J10").Value = 0.09058663713443794
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.972158333333333
$ws.Range("N10").Value = 23.916475
$ws.Range("O10").Value = 0.3253124204239203
$ws.Range("P10").Value = 0.3253124204239203
$ws.Range("Q10").Value = 258.3171375868111
$ws.Range("R10").Value = 2324.8542382813
$ws.Range("S10").Value = 0.02946895818426739
$ws.Range("T10").Value = 0.02946895818426739
